$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# "add consump to stock" - bump the stock/consumption (سرمایه) figure for the
# latest quarter (M26) and give the capital row (L26:M26) a thousands-separated
# "Comma" number format to match the rest of the statement.
$ws.Range("M26").Value = 7500000

$rng = $ws.Range("L26:M26")
$rng.NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"

# Reposition the view: scroll the window over to the later quarters, switch
# the sheet to right-to-left reading order, and leave the selection on M31.
$ws.Range("M31").Select()
$ws.Application.ActiveWindow.DisplayRightToLeft = $true
$ws.Application.ActiveWindow.ScrollColumn = 8
